$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New CUSTOMERS rows (TC, Customer_ID, PD) appended by this commit.
# Each entry: row number, then the text for columns A, B, C. All values in
# this sheet are stored as text/shared-strings (even the purely-numeric
# looking IDs), so every value must be force-written as text instead of
# letting Excel infer a number.
$newRows = @(
    @(92, "118463", "17705034", "1001"),
    @(93, "118463", "17705035", "1002"),
    @(94, "118463", "17705036", "1005"),
    @(95, "118464", "17705037", "1006"),
    @(96, "118465", "17705038", "1007"),
    @(97, "118466", "17705039", "1007"),
    @(98, "118468", "17705040", "1177"),
    @(99, "118469", "17705041", "1068"),
    @(100, "118470", "17705042", "6007"),
    @(101, "118471", "17705043", "1050"),
    @(102, "118463", "17705205", "1001"),
    @(103, "118463", "17705206", "1005"),
    @(104, "118464", "17705208", "1010"),
    @(105, "118466", "17705210", "6001"),
    @(106, "118468", "17705211", "6004"),
    @(107, "118469", "17705212", "6005"),
    @(108, "118470", "17705214", "1003"),
    @(109, "118471", "17705215", "1033"),
    @(110, "118463", "17705224", "1010"),
    @(111, "118463", "17705225", "1005"),
    @(112, "118463", "17705226", "1007"),
    @(113, "118463", "17705229", ""),
    @(114, "118463", "17705230", ""),
    @(115, "118463", "17705231", ""),
    @(116, "118464", "17705232", ""),
    @(117, "118465", "17705233", ""),
    @(118, "118466", "17705234", ""),
    @(119, "118468", "17705235", ""),
    @(120, "118471", "17705238", ""),
    @(121, "", "17704245", ""),
    @(122, "", "17704413", ""),
    @(123, "118463", "17705259", "1010"),
    @(124, "118463", "17705260", "1003"),
    @(125, "118463", "17705261", "1007"),
    @(126, "118464", "17705262", "1010"),
    @(127, "118465", "17705263", "1011")
)

foreach ($row in $newRows) {
    $r = $row[0]
    for ($col = 1; $col -le 3; $col++) {
        $text = $row[$col]
        $cell = $ws.Cells.Item($r, $col)
        if ($text -eq "") {
            # Force an explicit empty-text cell (t="s" pointing at the "" shared
            # string) rather than leaving the cell completely blank/absent.
            $cell.Value = "'"
        } else {
            # Force text storage so purely-numeric IDs do not turn into numbers.
            $cell.NumberFormat = "@"
            $cell.Value = $text
        }
    }
}
